$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.041.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "'2.645.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'581.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'156.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -2.80%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'2.643.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").Value = "'5.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'0.384"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "'28.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "'3.123.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "'63.932.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "'2.644.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'12.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "'7.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.97%  "
$ws.Range("E21").Value = "  -3.67%  "
$ws.Range("D22").Value = "'346.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "'68.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'1.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.02%  "
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").Value = "'9.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").Value = "'1.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").Value = "'587.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").Value = "'8.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "'1.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'6.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'5.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").Value = "'0.404"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("D38").Value = "'19.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").Value = "'151.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").Value = "'2.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.94%  "
$ws.Range("D44").Value = "'41.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "'163.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("D46").Value = "'24.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.11%  "
$ws.Range("D47").Value = "'3.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("D48").Value = "'0.0592"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").Value = "'0.635"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").Value = "  -2.53%  "
